$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename / retarget the header row (row 1)
# ------------------------------------------------------------------
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
# G1 currently holds "Q"; give it the new header text, but first copy the
# bold/centered/bordered look from A1 so the new header cell keeps the
# same styling the rest of row 1 already has.
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "has_multi_node_motifs"
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Drop the now-unused trailing columns H:L entirely (headers + data)
# ------------------------------------------------------------------
$ws.Range("H1:L2").Clear()

# ------------------------------------------------------------------
# 3) Insert three new data rows above the existing data row (old row 2),
#    pushing it down to row 5. Newly inserted rows inherit the header's
#    bold formatting by default, so strip that back off afterwards.
# ------------------------------------------------------------------
$ws.Range("A2:A4").EntireRow.Insert()
$ws.Range("B2:G4").ClearFormats()

# Give column A of every data row (2-5) the same styling as A1/old-A2 had.
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4) Populate the four data rows
# ------------------------------------------------------------------
$data = @(
    @{ Row = 2
       Glycan = "Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc"
       BindingScore = 0.4240903392243248
       Monosaccharides = "['Gal(b1-4)', 'GlcNAc(b1-3)']"
       Motifs = "['Gal(b1-4)GlcNAc(b1-3)']"
       Sasa = 5.338285572580087
       Flexibility = 0.91117855161729
       HasMultiNodeMotifs = $true },
    @{ Row = 3
       Glycan = "Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc"
       BindingScore = 1.73760896949211
       Monosaccharides = "['Gal(b1-4)', 'GlcNAc(b1-3)']"
       Motifs = "['Gal(b1-4)GlcNAc(b1-3)']"
       Sasa = 5.27278254643194
       Flexibility = 2.180924532303609
       HasMultiNodeMotifs = $true },
    @{ Row = 4
       Glycan = "Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc"
       BindingScore = 1.171840310090835
       Monosaccharides = "['Gal(b1-4)', 'GlcNAc(b1-3)']"
       Motifs = "['Gal(b1-4)GlcNAc(b1-3)']"
       Sasa = 5.252147263686476
       Flexibility = 1.839554809126105
       HasMultiNodeMotifs = $true },
    @{ Row = 5
       Glycan = "GlcNAc(b1-4)GlcNAc(b1-4)GlcNAc"
       BindingScore = 1.922476871100382
       Monosaccharides = "['GlcNAc(b1-4)', 'GlcNAc(b1-4)']"
       Motifs = "['GlcNAc(b1-4)GlcNAc']"
       Sasa = 5.702079978569953
       Flexibility = 0.6713220512263312
       HasMultiNodeMotifs = $true }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.Glycan
    $ws.Cells.Item($r, 2).Value = $row.BindingScore
    $ws.Cells.Item($r, 3).Value = $row.Monosaccharides
    $ws.Cells.Item($r, 4).Value = $row.Motifs
    $ws.Cells.Item($r, 5).Value = $row.Sasa
    $ws.Cells.Item($r, 6).Value = $row.Flexibility
    $ws.Cells.Item($r, 7).Value = $row.HasMultiNodeMotifs
}
